# ModelComponentClassDiagram.pptx - "Updating the diagrams in Developer Guide"
#
# 1) Rename class-diagram shapes from the AddressBook-ish domain to the
#    BookShelf-ish domain (AB3 -> AB4 sample app naming).
# 2) Resize the renamed "Address"->"ReviewList" box (and best-effort nudge
#    the elbow connector hanging off it) to match the new, wider label.
# 3) Refresh the cached text of the auto date placeholders (slide master +
#    every slide layout) the way PowerPoint itself re-stamps them on save.

$p = $ppt.ActivePresentation

# Replaces just the $oldSub substring inside the shape's text with $newSub,
# editing the existing run's Characters range in place. This avoids
# rebuilding the whole paragraph (which would flatten sibling runs / line
# breaks / fields), so untouched runs keep their original formatting.
function Replace-ShapeSubstring($shape, [string]$oldSub, [string]$newSub) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($oldSub)
    if ($idx -lt 0) { return }
    $tr.Characters($idx + 1, $oldSub.Length).Text = $newSub
}

# --- 1) Class / member name renames on the single slide ------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }

    switch ($sh.Id) {
        46  { Replace-ShapeSubstring $sh "VersionedAddressBook" "VersionedBookShelf" }
        49  { Replace-ShapeSubstring $sh "UniquePersonList" "UniqueBookList" }
        62  { Replace-ShapeSubstring $sh "Person" "Book" }
        80  { Replace-ShapeSubstring $sh "Phone" "Author" }
        83  { Replace-ShapeSubstring $sh "Email" "Rating" }
        85  { Replace-ShapeSubstring $sh "Address" "ReviewList" }
        100 { Replace-ShapeSubstring $sh "ReadOnlyAddressBook" "ReadOnlyBookShelf" }
        55  { Replace-ShapeSubstring $sh "AddressBook" "BookShelf" }
    }
}

# --- 2) Resize the "ReviewList" box (wider label than "Address") ---------
$reviewListBox = $null
$reviewConnector = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 85) { $reviewListBox = $sh }
    if ($sh.Id -eq 86) { $reviewConnector = $sh }
}
if ($reviewListBox -ne $null) {
    $reviewListBox.Left = 607.2753
    $reviewListBox.Width = 63.97787401574803
}
if ($reviewConnector -ne $null) {
    $reviewConnector.Width = 34.204803
}

# --- 3) Re-cache the "today" date placeholders (master + every layout) ---
function Update-DatePlaceholders($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if (-not $sh.HasTextFrame) { continue }
        if ($sh.TextFrame.TextRange.Text -eq "12/5/2018") {
            Replace-ShapeSubstring $sh "12/5/2018" "3/4/2019"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholders $layout.Shapes
}

Write-Host "Edit applied."
